$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new header cells (left to right so shared-string table is
#     built in the same order as in the target workbook) ---
$ws.Range("C1").Value = "Hour"
$ws.Range("D1").Value = "air_temp_Cogealac"
$ws.Range("E1").Value = "wind_direction_100m_Cogealac"
$ws.Range("F1").Value = "wind_direction_10m_Cogealac"
$ws.Range("G1").Value = "wind_speed_100m_Cogealac"
$ws.Range("H1").Value = "wind_speed_10m_Cogealac"
$ws.Range("I1").Value = "wind_gust_Cogealac"
$ws.Range("J1").Value = "air_temp_Focsani"
$ws.Range("K1").Value = "wind_direction_100m_Focsani"
$ws.Range("L1").Value = "wind_direction_10m_Focsani"
$ws.Range("M1").Value = "wind_speed_100m_Focsani"
$ws.Range("N1").Value = "wind_speed_10m_Focsani"
$ws.Range("O1").Value = "wind_gust_Focsani"

# --- Apply the "Normal 2" header style (same as B1/Interval) to the new
#     Hour, air_temp_Cogealac and air_temp_Focsani header cells ---
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("B1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("B1").Copy()
$ws.Range("J1").PasteSpecial(-4122)

# --- Resize the columns to fit their new headers ---
$ws.Columns.Item(3).ColumnWidth = 4.0
$ws.Columns.Item(4).ColumnWidth = 16.0
$ws.Columns.Item(5).ColumnWidth = 26.0
$ws.Columns.Item(6).ColumnWidth = 25.0
$ws.Columns.Item(7).ColumnWidth = 23.68
$ws.Columns.Item(8).ColumnWidth = 22.68
$ws.Columns.Item(9).ColumnWidth = 16.68
$ws.Columns.Item(10).ColumnWidth = 14.68
$ws.Columns.Item(11).ColumnWidth = 24.68
$ws.Columns.Item(12).ColumnWidth = 23.68
$ws.Columns.Item(13).ColumnWidth = 22.34
$ws.Columns.Item(14).ColumnWidth = 21.34
$ws.Columns.Item(15).ColumnWidth = 15.34

# --- Match the saved selection state ---
$ws.Range("F14").Select() | Out-Null
